$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets   = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
# Core semantic change: the Orchestrator queue folder now includes the
# NewHireCommunication sub-folder (building name included downstream).
$wsSettings.Range("B3").Value = "DEV/P004_NewHireCommunication"

# OrchestratorFolder (row 21) simply mirrors OrchestratorQueueFolder (B3)
# rather than duplicating the literal value.
$wsSettings.Range("B21").Formula = "=B3"

# Three new Name/Value rows referencing the shared O365 asset names that used
# to live only on the Assets sheet.
$wsSettings.Range("A25").Value = "O365TenantID"
$wsSettings.Range("B25").Value = "Shared_O365TenantID"
$wsSettings.Range("A26").Value = "O365AppID"
$wsSettings.Range("B26").Value = "Shared_O365ApplicationID"
$wsSettings.Range("A27").Value = "O365ApplicationSecret"
$wsSettings.Range("B27").Value = "Shared_O365ApplicationSecret"

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
# The shared O365 app/tenant/secret assets are no longer needed here - they
# are now referenced from the Settings sheet (rows above) instead.
$wsAssets.Rows("4:6").Delete()

# The OrchestratorAssetFolder column (C) is now computed from the Settings
# sheet folder value instead of being hard-coded, so that both sheets always
# agree. Rows that belong to the SP_003_WorkdayDisposition sub-process append
# their own sub-folder.
$plainRows = @(2,3,4,5,6,7,8,12,18,23,24,25)
$subRows   = @(9,10,11,13,14,15,16,17,19,20,21,22)

foreach ($r in $plainRows) {
    $wsAssets.Range("C$r").Formula = "=Settings!B3"
}
foreach ($r in $subRows) {
    $wsAssets.Range("C$r").Formula = "=Settings!B3&""/SP_003_WorkdayDisposition"""
}

# ---------------------------------------------------------------------------
# Workbook level bookkeeping
# ---------------------------------------------------------------------------
# The hidden _FilterDatabase name on Assets shrank along with the data range
# after the three asset rows were removed.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Assets!_FilterDatabase") {
        $n.RefersTo = "=Assets!`$A`$1:`$D`$25"
    }
}

# Restore the selections/active sheet as left by the author.
$wsAssets.Range("C30").Select()
$wsSettings.Activate()
$wsSettings.Range("B25").Select()
